# Bug 2920 Updated docs
#
# Duplicate the active "V1.2.2 7_10_2015" sheet, insert the copy before it,
# rename the copy to "V1.2.3 7_31_201" and mark a few of its checklist
# checkmarks as not-yet-done (clearing some "x" cells + updating one
# instruction cell), since this new sheet tracks the NEXT release.

$wb = $excel.ActiveWorkbook

# The sheet that is currently first/active ("V1.2.2 7_10_2015").
$original = $wb.Worksheets.Item(1)

# Duplicate it and place the new copy immediately before the original -
# this becomes the new first/active tab.
$original.Copy($original)

$newSheet = $wb.Worksheets.Item(1)
$newSheet.Name = "V1.2.3 7_31_201"

# Clear a handful of "done" checkmarks on the new release's sheet - these
# tests have not been re-run yet for V1.2.3.
$newSheet.Range("E6:E9").ClearContents()
$newSheet.Range("D10:F13").ClearContents()

# Update the release-directory verification instructions for the new
# version number.
$newSheet.Range("J18").Value = "Confirm release directory has necessary files + ex: /Users/scoleman/dev/fips/fcids/release/V1.2.0_01_26_2015 +Confirm no errors"

# Reflect the user's on-screen selection at save time: the new sheet has
# E6:E9 selected (the range just cleared) and is the active tab, while the
# duplicated-from sheet keeps a separate, different selection.
$newSheet.Range("E6:E9").Select()

$original.Select()
$original.Range("B32").Select()

$newSheet.Select()
